$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Splněno?" answer for the task in row 22 ("Formátování tabulek pro
# výpis z databáze") changes from "ne" to "ano" - the underlying dependent
# formulas (F22, and the J4/J6/J7 summary counters) recalc automatically.
$ws.Range("E22").Value = "ano"

# Move the active selection to match where the user clicked next.
$ws.Range("G22").Select()
